$wb = $excel.ActiveWorkbook

# --- A1 sheet: add submission row (TxHash, ClassID) ---
$ws = $wb.Worksheets.Item("A1")
$ws.Activate() | Out-Null
$ws.Range("A3").Value = "D4677935A5E55F7EBD33E02998F5B9E312BDF8D71CB803DF39279142D31CBD3F"
$ws.Range("B3").Value = "conqueror"
$ws.Range("B6").Select() | Out-Null

# --- A2 sheet: add two submission rows (TxHash, ClassID, NFTID) ---
$ws = $wb.Worksheets.Item("A2")
$ws.Activate() | Out-Null
$ws.Range("A3").Value = "4DB15DC325E6E2B3C8E818C64B085320A7E553621533DB554BBA60EFEE71BBD1"
$ws.Range("B3").Value = "conqueror"
$ws.Range("C3").Value = "conquerornft1"
$ws.Range("A4").Value = "D10102A3DF7FD0256B18041DD948543644E78DFFF7488EFCF3A7919A38F6C2FF"
$ws.Range("B4").Value = "conqueror"
$ws.Range("C4").Value = "conquerornft2"
$ws.Range("C4").Select() | Out-Null

# --- A3 sheet: Juno destination submission (TxHash, ClassID, NFTID, ChainID) ---
$ws = $wb.Worksheets.Item("A3")
$ws.Activate() | Out-Null
$ws.Range("A3").Value = "CAC9A9AA9A6D572C67B18F29245FC5E4C06D9D09AA6C46AFA3CAF92DCEC23C7A"
$ws.Range("B3").Value = "juno1jm3qhjnwt8qr0sjcmhg9rp4zr79twyq5h0wupf0fffmr6rrztqvq7fn7cf"
$ws.Range("C3").Value = "conquerornft1"
$ws.Range("D3").Value = "uni-6"
$ws.Range("D3").Select() | Out-Null

# --- A4 sheet: OmniFlix destination submission (TxHash, ClassID, NFTID, ChainID) ---
$ws = $wb.Worksheets.Item("A4")
$ws.Activate() | Out-Null
$ws.Range("A3").Value = "C64A644724A0C6849DB8F38DD9D8467CCE14D6762B52F4F3AF033D6268244E04"
$ws.Range("B3").Value = "ibc/ADF39F749EE1927F3518A2EFA47FA153AD553167CBE3AF40ACD38D75D3137E83"
$ws.Range("C3").Value = "conquerornft2"
$ws.Range("D3").Value = "gon-flixnet-1"
$ws.Range("D3").Select() | Out-Null

# --- A5 sheet: Juno destination submission again (internal transfer row) ---
$ws = $wb.Worksheets.Item("A5")
$ws.Activate() | Out-Null
$ws.Range("A3").Value = "57F37DA3ADCE38D0F6C5D89FCCA2AB414322894ADDB07986CBF0A315106D497D"
$ws.Range("B3").Value = "juno1jm3qhjnwt8qr0sjcmhg9rp4zr79twyq5h0wupf0fffmr6rrztqvq7fn7cf"
$ws.Range("C3").Value = "conquerornft1"
$ws.Range("D3").Value = "uni-6"
$ws.Range("D3").Select() | Out-Null

# --- A6 sheet: OmniFlix destination submission again (internal transfer row) ---
$ws = $wb.Worksheets.Item("A6")
$ws.Activate() | Out-Null
$ws.Range("A3").Value = "3F88138C5904C4D9BBC4E9EA7DC2B75CD536F4ECB3690EA64828F0BC0DAE8446"
$ws.Range("B3").Value = "ibc/ADF39F749EE1927F3518A2EFA47FA153AD553167CBE3AF40ACD38D75D3137E83"
$ws.Range("C3").Value = "conquerornft2"
$ws.Range("D3").Value = "gon-flixnet-1"
$ws.Range("D3").Select() | Out-Null
